# October 2014 payslips - "pattabhi - modifieed schedulers"
# Remove the unused/sample helper columns (W:AC) so the sheet collapses
# from A1:AC5 down to A1:V5, then refresh the payroll columns (N:V) and a
# handful of per-employee fields with the updated schedule data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the trailing placeholder/duplicate columns (old W:AC) -----------
$ws.Columns("W:AC").Delete()

# --- Column widths (character units); ColumnWidth = stored width - 5/7 ----
function Set-ColWidth($colLetter, $storedWidth) {
    $ws.Columns($colLetter + ":" + $colLetter).ColumnWidth = $storedWidth - (5.0 / 7.0)
}

Set-ColWidth "C" 15.18988764044944
Set-ColWidth "F" 22.88988764044944
Set-ColWidth "G" 18.48988764044944
Set-ColWidth "N" 11.88988764044944
Set-ColWidth "O" 5.289887640449439
Set-ColWidth "P" 9.68988764044944
Set-ColWidth "Q" 5.289887640449439
Set-ColWidth "R" 6.389887640449439
Set-ColWidth "S" 5.289887640449439
Set-ColWidth "T" 6.389887640449439
Set-ColWidth "U" 11.88988764044944
Set-ColWidth "V" 9.68988764044944

# --- Header row (payroll columns) ------------------------------------------
$ws.Range("N1").Value = "Spcl Allowance"
$ws.Range("O1").Value = "Arrears"
$ws.Range("P1").Value = "Gross Pay"
$ws.Range("Q1").Value = "PF"
$ws.Range("R1").Value = "ESIC"
$ws.Range("S1").Value = "PT"
$ws.Range("T1").Value = "TDS"
$ws.Range("U1").Value = "total_deducations"
$ws.Range("V1").Value = "NetPay"

# --- Row 2: Vidya Sagar  Pogiri (was Sekhar Beri) ---------------------------
$ws.Range("C2").Value = "Vidya Sagar  Pogiri"
$ws.Range("F2").Value = "Junior Development"
$ws.Range("G2").Value = "Development"
$ws.Range("J2").Value = 10000.0
$ws.Range("N2").Value = 6000.0
$ws.Range("O2").Value = 0.0
$ws.Range("P2").Value = 10000.0
$ws.Range("Q2").Value = 0.0
$ws.Range("R2").Value = 0.0
$ws.Range("S2").Value = 0.0
$ws.Range("T2").Value = 0.0
$ws.Range("U2").Value = 0.0
$ws.Range("V2").Value = 10000.0

# --- Row 3: BalaRaju Vankala -------------------------------------------------
$ws.Range("F3").Value = "Junior Accounts"
$ws.Range("G3").Value = "Accounts"
$ws.Range("J3").Value = 10000.0
$ws.Range("N3").Value = 6000.0
$ws.Range("O3").Value = 0.0
$ws.Range("P3").Value = 10000.0
$ws.Range("Q3").Value = 0.0
$ws.Range("R3").Value = 0.0
$ws.Range("S3").Value = 0.0
$ws.Range("T3").Value = 0.0
$ws.Range("U3").Value = 0.0
$ws.Range("V3").Value = 10000.0

# --- Row 4: Priyanka Muddana -------------------------------------------------
$ws.Range("F4").Value = "Junior Business Development"
$ws.Range("G4").Value = "Business Development"
$ws.Range("J4").Value = 10000.0
$ws.Range("N4").Value = 6000.0
$ws.Range("O4").Value = 0.0
$ws.Range("P4").Value = 10000.0
$ws.Range("Q4").Value = 0.0
$ws.Range("R4").Value = 0.0
$ws.Range("S4").Value = 0.0
$ws.Range("T4").Value = 0.0
$ws.Range("U4").Value = 0.0
$ws.Range("V4").Value = 10000.0

# --- Row 5: pattabhi ramarao --------------------------------------------------
$ws.Range("F5").Value = "Junior HR"
$ws.Range("G5").Value = "HR"
$ws.Range("J5").Value = 10000.0
$ws.Range("N5").Value = 6000.0
$ws.Range("O5").Value = 0.0
$ws.Range("P5").Value = 10000.0
$ws.Range("Q5").Value = 0.0
$ws.Range("R5").Value = 0.0
$ws.Range("S5").Value = 0.0
$ws.Range("T5").Value = 0.0
$ws.Range("U5").Value = 0.0
$ws.Range("V5").Value = 10000.0
